# Updates cryptos list values/links per the Dec 13 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    # Force text storage for values that would otherwise be auto-parsed as a
    # number by Excel (e.g. "245.57"), then restore the General format so the
    # cell keeps its original (unstyled) number format.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

$ws.Range("D2").Value = "40.830.56"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.158.63"
$ws.Range("E3").Value = "  -3.26%  "
$ws.Range("E4").Value = "  -0.23%  "
Set-TextCell "D5" "245.57"
$ws.Range("E5").Value = "  -2.43%  "
Set-TextCell "D6" "0.612"
$ws.Range("E6").Value = "  -2.80%  "
Set-TextCell "D7" "65.75"
$ws.Range("E7").Value = "  -7.60%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.03%  "
Set-TextCell "D10" "57.85"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D11" "35.45"
$ws.Range("E11").Value = "  -15.26%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D12" "0.0917"
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("E13").Value = "  -1.94%  "
Set-TextCell "D14" "6.84"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "2.474.55"
$ws.Range("E15").Value = "  -3.27%  "
Set-TextCell "D16" "0.849"
$ws.Range("E16").Value = "  -0.93%  "
Set-TextCell "D17" "14.14"
$ws.Range("E17").Value = "  -5.49%  "
$ws.Range("D18").Value = "2.172.82"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").Value = "40.684.35"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "0.0₃0930"
$ws.Range("E20").Value = "  -3.90%  "
Set-TextCell "D21" "6.04"
$ws.Range("E21").Value = "  -2.02%  "
Set-TextCell "D22" "70.85"
$ws.Range("E22").Value = "  -2.93%  "
Set-TextCell "D23" "227.24"
$ws.Range("E23").Value = "  -3.08%  "
Set-TextCell "D24" "2.08"
$ws.Range("E24").Value = "  -7.24%  "
Set-TextCell "D25" "11.47"
$ws.Range("E25").Value = "  +12.65%  "
$ws.Range("E26").Value = "  -0.01%  "
Set-TextCell "D27" "3.65"
$ws.Range("E27").Value = "  -2.57%  "
Set-TextCell "D28" "2.38"
$ws.Range("E28").Value = "  -5.36%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D29" "2.12"
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D30" "167.45"
$ws.Range("E30").Value = "  -1.60%  "
Set-TextCell "D31" "20.02"
$ws.Range("E31").Value = "  -3.01%  "
Set-TextCell "D32" "0.118"
$ws.Range("E32").Value = "  -1.83%  "
Set-TextCell "D33" "5.53"
$ws.Range("E33").Value = "  +0.73%  "
Set-TextCell "D34" "0.0728"
$ws.Range("E34").Value = "  +1.10%  "
Set-TextCell "D35" "0.120"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D36" "4.53"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D37" "24.93"
$ws.Range("E37").Value = "  -6.19%  "
Set-TextCell "D38" "3.94"
$ws.Range("E38").Value = "  -3.56%  "
Set-TextCell "D39" "0.0295"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("E40").Value = "  -5.42%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell "D41" "5.40"
$ws.Range("E41").Value = "  -10.12%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D42" "11.41"
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell "D43" "60.06"
$ws.Range("E43").Value = "  -14.91%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D44" "4.78"
$ws.Range("E44").Value = "  -7.01%  "
$ws.Range("E45").Value = "  -10.95%  "
$ws.Range("E46").Value = "  -0.08%  "
Set-TextCell "D47" "8.38"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("E51").Value = "  -1.55%  "
